{"js": "// Update the date line and all \"two-digit \u00f7 one-digit\" problems in the\n// table to the new values from the day's worksheet, in document order.\n// (Source: 2024-01-31 Wednesday -> 2024-02-01 Thursday worksheet refresh.)\nconst replacements = [\n  [\"2024-01-31 Wednesday\", \"2024-02-01 Thursday\"],\n  [\"13\u00f75=\", \"51\u00f78=\"],\n  [\"94\u00f72=\", \"74\u00f75=\"],\n  [\"63\u00f73=\", \"16\u00f78=\"],\n  [\"72\u00f78=\", \"18\u00f73=\"],\n  [\"79\u00f73=\", \"27\u00f78=\"],\n  [\"90\u00f73=\", \"84\u00f79=\"],\n  [\"91\u00f77=\", \"94\u00f79=\"],\n  [\"73\u00f77=\", \"57\u00f79=\"],\n  [\"63\u00f74=\", \"88\u00f79=\"],\n  [\"59\u00f78=\", \"81\u00f79=\"],\n  [\"34\u00f78=\", \"64\u00f79=\"],\n  [\"41\u00f72=\", \"66\u00f77=\"],\n  [\"24\u00f76=\", \"72\u00f78=\"],\n  [\"35\u00f75=\", \"44\u00f73=\"],\n  [\"33\u00f76=\", \"29\u00f73=\"],\n  [\"40\u00f74=\", \"16\u00f72=\"],\n  [\"31\u00f74=\", \"52\u00f75=\"],\n  [\"47\u00f75=\", \"51\u00f77=\"],\n  [\"40\u00f74=\", \"31\u00f79=\"],\n  [\"92\u00f79=\", \"53\u00f76=\"],\n  [\"83\u00f74=\", \"95\u00f73=\"],\n  [\"65\u00f75=\", \"43\u00f72=\"],\n  [\"85\u00f79=\", \"79\u00f72=\"],\n  [\"77\u00f79=\", \"66\u00f73=\"],\n  [\"10\u00f73=\", \"73\u00f75=\"],\n];\n\nconst body = context.document.body;\n\n// Cache one search-results collection per distinct \"old\" string, and pop\n// occurrences off the front in document order. This correctly handles the\n// \"40\u00f74=\" value, which appears twice in the source with two different\n// targets.\nconst resultsCache = new Map();\nconst consumedCount = new Map();\n\nfor (const [oldText] of replacements) {\n  if (!resultsCache.has(oldText)) {\n    const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    found.load(\"items\");\n    resultsCache.set(oldText, found);\n    consumedCount.set(oldText, 0);\n  }\n}\nawait context.sync();\n\nfor (const [oldText, newText] of replacements) {\n  const found = resultsCache.get(oldText);\n  const n = consumedCount.get(oldText);\n  const range = found.items[n];\n  range.insertText(newText, Word.InsertLocation.replace);\n  consumedCount.set(oldText, n + 1);\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and all \"two-digit \u00f7 one-digit\" problems in the\n# table to the new values from the day's worksheet, in document order.\n# (Source: 2024-01-31 Wednesday -> 2024-02-01 Thursday worksheet refresh.)\n$pairs = @(\n    @(\"2024-01-31 Wednesday\", \"2024-02-01 Thursday\"),\n    @(\"13\u00f75=\", \"51\u00f78=\"),\n    @(\"94\u00f72=\", \"74\u00f75=\"),\n    @(\"63\u00f73=\", \"16\u00f78=\"),\n    @(\"72\u00f78=\", \"18\u00f73=\"),\n    @(\"79\u00f73=\", \"27\u00f78=\"),\n    @(\"90\u00f73=\", \"84\u00f79=\"),\n    @(\"91\u00f77=\", \"94\u00f79=\"),\n    @(\"73\u00f77=\", \"57\u00f79=\"),\n    @(\"63\u00f74=\", \"88\u00f79=\"),\n    @(\"59\u00f78=\", \"81\u00f79=\"),\n    @(\"34\u00f78=\", \"64\u00f79=\"),\n    @(\"41\u00f72=\", \"66\u00f77=\"),\n    @(\"24\u00f76=\", \"72\u00f78=\"),\n    @(\"35\u00f75=\", \"44\u00f73=\"),\n    @(\"33\u00f76=\", \"29\u00f73=\"),\n    @(\"40\u00f74=\", \"16\u00f72=\"),\n    @(\"31\u00f74=\", \"52\u00f75=\"),\n    @(\"47\u00f75=\", \"51\u00f77=\"),\n    @(\"40\u00f74=\", \"31\u00f79=\"),\n    @(\"92\u00f79=\", \"53\u00f76=\"),\n    @(\"83\u00f74=\", \"95\u00f73=\"),\n    @(\"65\u00f75=\", \"43\u00f72=\"),\n    @(\"85\u00f79=\", \"79\u00f72=\"),\n    @(\"77\u00f79=\", \"66\u00f73=\"),\n    @(\"10\u00f73=\", \"73\u00f75=\")\n)\n\n$d = $word.ActiveDocument\n\n# Re-fetch $d.Content for every search and use wdReplaceOne (1) so each call\n# only consumes a single occurrence, in document order. This is required\n# because \"40\u00f74=\" appears twice in the source with two different targets;\n# doing them one at a time (instead of a global wdReplaceAll) keeps the\n# first occurrence mapped to \"16\u00f72=\" and the second to \"31\u00f79=\".\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 1) | Out-Null\n}\n"}
